$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a brand-new row above row 5 ("History of partnering with Wannon Water").
#    Use Insert() then copy formatting from the row that ends up right below it (new
#    row 6, which carries the style pattern we want: A=s6, B=s10, C=s17).
$ws.Rows(5).Insert()

$ws.Range("A6:C6").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A5").Value = "History of partnering with Wannon Water"
$ws.Range("C5").Value = "Very important!"

# 2) Extend the "Proposed deliverables" sub-criteria text (now row 12, column C) with
#    the extra sentence about calling out seats/licences.
$ws.Range("C12").Value = "The service will be delivered through the following channels:`nResearch Reports: Access to detailed reports on trends, best practices, and technologies in IT and OT.`nMarket Analysis: Regular updates on the competitive landscape, including vendor evaluations, market forecasts, and technology adoption rates.`nStrategic Guidance: Customised strategic advice based on the latest industry trends and specific business needs.`nBenchmarking Services: Data and tools to compare Wannon Water's performance against industry peers.`nAccess to Analysts: Direct consultations with industry experts for personalised advice and support.`nWorkshops, Webinars, and Conferences: Access to world-class conferences, educational sessions, and workshops that provide cutting-edge insights and networking opportunities with industry leaders.`nSpecifically call out the number of seats or employee access licences as part of the proposal."
$ws.Rows(12).RowHeight = 331.2

# 3) Defined names need to point at the criteria cells' new row numbers.
$wb.Names.Item("Sheet1!_Toc176179427").RefersTo = "=Sheet1!`$C`$16"
$wb.Names.Item("Sheet1!_Toc176179431").RefersTo = "=Sheet1!`$C`$13"
$wb.Names.Item("Sheet1!_Toc176179432").RefersTo = "=Sheet1!`$C`$17"
$wb.Names.Item("Sheet1!_Toc176179433").RefersTo = "=Sheet1!`$C`$20"

# 4) The hyperlink on the eServices contract cell moves from C25 to C26 (row shifted
#    by the inserted row). Re-create it in place, preserving the cell's own text.
$hyperlinkUrl = "https://content.vic.gov.au/sites/default/files/2023-12/eServcies-contract-%28April-2021%29.pdf"
$hyperlinkText = $ws.Range("C26").Text
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C26"), $hyperlinkUrl, "", "", $hyperlinkText)

# 5) Update the saved selection/scroll position to match the edited area.
$ws.Range("C12").Select()
